$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2 and Row 3 share identical target values per the diff
foreach ($r in 2,3) {
    $ws.Range("D" + $r).Value = 0.142
    $ws.Range("E" + $r).Value = -0.158
    $ws.Range("G" + $r).Value = 0.3081383201002925
    $ws.Range("H" + $r).Value = 0.3081383201002925
    $ws.Range("I" + $r).Value = 0.06963017133305474
    $ws.Range("J" + $r).Value = 0.06963017133305474
    $ws.Range("K" + $r).Value = 54.6
    $ws.Range("L" + $r).Value = 0.02852068533221897
    $ws.Range("M" + $r).Value = 83.5
    $ws.Range("N" + $r).Value = 0.04222930258433217
    $ws.Range("O" + $r).Value = 1.529304029304029
    $ws.Range("P" + $r).Value = 35.3
    $ws.Range("Q" + $r).Value = 0.01785262732008294
    $ws.Range("R" + $r).Value = 0.6465201465201464
    $ws.Range("S" + $r).Value = 48.2
    $ws.Range("T" + $r).Value = 0.577245508982036
    $ws.Range("U" + $r).Value = 5.66
    $ws.Range("V" + $r).Value = 0.002862489253021798
    $ws.Range("W" + $r).Value = 0.053377651774367
    $ws.Range("X" + $r).Value = 0.1048847983307459
    $ws.Range("Y" + $r).Value = -0.05150714655637888
    $ws.Range("Z" + $r).Value = 1.891662220114227
    $ws.Range("AA" + $r).Value = 0.1317167644908203
    $ws.Range("AB" + $r).Value = 0.1048218937552831
    $ws.Range("AC" + $r).Value = 0.02689487073553729
    $ws.Range("AD" + $r).Value = 1.77
    $ws.Range("AE" + $r).Value = 0
    $ws.Range("AF" + $r).Value = 1.77
    $ws.Range("AG" + $r).Value = -3.89
    $ws.Range("AH" + $r).Value = 0.0008943594718731526
    $ws.Range("AI" + $r).Value = 0.002035254751802408
    $ws.Range("AJ" + $r).Value = -0.001971207199720282
    $ws.Range("AK" + $r).Value = -0.004502262705292763
    $ws.Range("AL" + $r).Value = 130.5
    $ws.Range("AM" + $r).Value = 130.5
    $ws.Range("AN" + $r).Value = 0.01314031180400891
    $ws.Range("AO" + $r).Value = 1.021455938697318
    $ws.Range("AP" + $r).Value = -0.02887899034892354
    $ws.Range("AQ" + $r).Value = 1.021455938697318
    # Column F (expected_growth_eps_next_5_years) is cleared - no longer populated
    $ws.Range("F" + $r).ClearContents()
}

Write-Output "Done updating brazil_reinsurance rows 2 and 3"
